$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I; this shifts existing columns I..Q
# (and their widths / contents) one position to the right, to J..R.
$ws.Columns("I:I").Insert()

# New header cell for the inserted column.
$ws.Range("I1").Value = "RequestContent"
$ws.Range("B1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats -> default/Normal style (s=0)
$ws.Application.CutCopyMode = $false

# Row 5 gains the DDL "create table" request content plus a couple of
# supporting cells (matching the style already used elsewhere in the sheet).
$ws.Range("I5").Value = "create table employees (emp_no int, birth_date date,first_name VARCHAR(50),last_name VARCHAR(50), gender VARCHAR(50),hire_date date)"
$ws.Range("B5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("K3").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").ClearContents()
$ws.Application.CutCopyMode = $false

$ws.Range("O5").Value = "DDL"
$ws.Range("B5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New row 6: a second DB test case (INSERT statement).
$ws.Range("A6").Value = "PetPost"
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B6").Value = "DB"
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("C6").Value = "pet"
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("D6").Value = "pet "
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I6").Value = "insert into employees (emp_no,birth_date,first_name,last_name, gender,hire_date) values  (2,'1978-01-08','ELan', 'Thangamani', 'Male', '2007-10-10')"
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("N6").Value = "select * from employees where emp_no = 1 `nEMP_NO,BIRTH_DATE,FIRST_NAME,LAST_NAME,GENDER,HIRE_DATE`ni~1,l~253087200000,ELan,Thangamani,Male,l~1191992400000          "
$ws.Range("N5").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("O6").Value = "INSERT"
$ws.Range("O5").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fix the custom row height (setting the multi-line N6 value above made
# Excel auto-expand the row); restore it to match the committed file.
$ws.Rows("6:6").RowHeight = 11.45

# Selection ends on I2, matching the committed file.
[void]$ws.Range("I2").Select()
